$wb = $excel.ActiveWorkbook

# --- Sheet "Under 18" (Scotland pre-1994 predicted rates added) ---
$ws = $wb.Worksheets.Item("Under 18")

# Insert two new columns before column B, shifting 1987-2018 data from B:AG to D:AI
$ws.Range("B1:C1").EntireColumn.Insert()

# New year headers (literal values, matching how the original B1/C1 were hardcoded)
$ws.Range("B1").Value = 1985
$ws.Range("C1").Value = 1986

# New predicted Scotland conception-rate values for 1985-1993 (row 2), formatted to 1 decimal place
$ws.Range("B2").Value = 38.418366115568197
$ws.Range("C2").Value = 39.256642256568902
$ws.Range("D2").Value = 38.793512380761399
$ws.Range("E2").Value = 39.594226250004397
$ws.Range("F2").Value = 40.928227491229798
$ws.Range("G2").Value = 42.785383249327403
$ws.Range("H2").Value = 44.409212627182903
$ws.Range("I2").Value = 43.1995217721907
$ws.Range("J2").Value = 41.278529201879202
$ws.Range("B2:J2").NumberFormat = "0.0"

# Make this the active sheet / selection, matching the authored workbook view
$ws.Range("B2:J2").Select()
$excel.ActiveWindow.ActiveSheet.Name | Out-Null
$wb.Worksheets.Item("Sheet1").Select()
$ws.Select()
